$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.039483
$ws.Range("H2").Value = 3.118449
$ws.Range("I2").Value = 0.01443400247257094
$ws.Range("J2").Value = 0.01491631751620889
$ws.Range("M2").Value = 8.850733666666665
$ws.Range("N2").Value = 26.552201
$ws.Range("O2").Value = 0.03425243510433672
$ws.Range("P2").Value = 0.03471479367970558
$ws.Range("Q2").Value = 9.200187184027664
$ws.Range("R2").Value = 82.80168465624899
$ws.Range("S2").Value = 0.000494399732987572
$ws.Range("T2").Value = 0.0005178168850361699
$ws.Range("G3").Value = 1.039483
$ws.Range("H3").Value = 3.118449
$ws.Range("I3").Value = 0.01443400247257094
$ws.Range("J3").Value = 0.01491631751620889
$ws.Range("O3").Value = 0.06418224226698629
$ws.Range("P3").Value = 0.06504861016194469
$ws.Range("Q3").Value = 17.23931863379033
$ws.Range("R3").Value = 155.153867704113
$ws.Range("S3").Value = 0.0009264066435768275
$ws.Range("T3").Value = 0.0009702857231636591
$ws.Range("G4").Value = 1.039483
$ws.Range("H4").Value = 3.118449
$ws.Range("I4").Value = 0.01443400247257094
$ws.Range("J4").Value = 0.01491631751620889
$ws.Range("M4").Value = 94.82067633333332
$ws.Range("N4").Value = 284.462029
$ws.Range("O4").Value = 0.3669570438989389
$ws.Range("P4").Value = 0.3719104358409092
$ws.Range("Q4").Value = 98.56448109700231
$ws.Range("R4").Value = 887.0803298730209
$ws.Range("S4").Value = 0.005296658878964609
$ws.Range("T4").Value = 0.005547534148594636
$ws.Range("G5").Value = 1.039483
$ws.Range("H5").Value = 3.118449
$ws.Range("I5").Value = 0.01443400247257094
$ws.Range("J5").Value = 0.01491631751620889
$ws.Range("M5").Value = 10.324604
$ws.Range("N5").Value = 20.649208
$ws.Range("O5").Value = 0.03995632924983982
$ws.Range("P5").Value = 0.02699712145781534
$ws.Range("Q5").Value = 10.732250339732
$ws.Range("R5").Value = 64.393502038392
$ws.Range("S5").Value = 0.0005767297551870467
$ws.Range("T5").Value = 0.0004026976356884298
$ws.Range("G6").Value = 1.039483
$ws.Range("H6").Value = 3.118449
$ws.Range("I6").Value = 0.01443400247257094
$ws.Range("J6").Value = 0.01491631751620889
$ws.Range("M6").Value = 127.8166836666667
$ws.Range("N6").Value = 383.450051
$ws.Range("O6").Value = 0.4946519494798983
$ws.Range("P6").Value = 0.5013290388596253
$ws.Range("Q6").Value = 132.8632697878776
$ws.Range("R6").Value = 1195.769428090899
$ws.Range("S6").Value = 0.00713980746185489
$ws.Range("T6").Value = 0.007477983123725996
$ws.Range("I7").Value = 0.8791289547788569
$ws.Range("J7").Value = 0.9085052224491242
$ws.Range("M7").Value = 8.850733666666665
$ws.Range("N7").Value = 26.552201
$ws.Range("O7").Value = 0.03425243510433672
$ws.Range("P7").Value = 0.03471479367970558
$ws.Range("Q7").Value = 560.3539945509956
$ws.Range("R7").Value = 5043.18595095896
$ws.Range("S7").Value = 0.03011230747190616
$ws.Range("T7").Value = 0.03153857135425637
$ws.Range("I8").Value = 0.8791289547788569
$ws.Range("J8").Value = 0.9085052224491242
$ws.Range("O8").Value = 0.06418224226698629
$ws.Range("P8").Value = 0.06504861016194469
$ws.Range("S8").Value = 0.05642446755953903
$ws.Range("T8").Value = 0.05909700204518392
$ws.Range("I9").Value = 0.8791289547788569
$ws.Range("J9").Value = 0.9085052224491242
$ws.Range("M9").Value = 94.82067633333332
$ws.Range("N9").Value = 284.462029
$ws.Range("O9").Value = 0.3669570438989389
$ws.Range("P9").Value = 0.3719104358409092
$ws.Range("Q9").Value = 6003.247499076674
$ws.Range("R9").Value = 54029.22749169006
$ws.Range("S9").Value = 0.3226025624516133
$ws.Range("T9").Value = 0.3378825732447959
$ws.Range("I10").Value = 0.8791289547788569
$ws.Range("J10").Value = 0.9085052224491242
$ws.Range("M10").Value = 10.324604
$ws.Range("N10").Value = 20.649208
$ws.Range("O10").Value = 0.03995632924983982
$ws.Range("P10").Value = 0.02699712145781534
$ws.Range("Q10").Value = 653.6670643865481
$ws.Range("R10").Value = 3922.002386319289
$ws.Range("S10").Value = 0.03512676597021155
$ws.Range("T10").Value = 0.02452702583551855
$ws.Range("I11").Value = 0.8791289547788569
$ws.Range("J11").Value = 0.9085052224491242
$ws.Range("M11").Value = 127.8166836666667
$ws.Range("N11").Value = 383.450051
$ws.Range("O11").Value = 0.4946519494798983
$ws.Range("P11").Value = 0.5013290388596253
$ws.Range("Q11").Value = 8092.277088013646
$ws.Range("R11").Value = 72830.49379212281
$ws.Range("S11").Value = 0.4348628513255869
$ws.Range("T11").Value = 0.4554600499693696
$ws.Range("G12").Value = 0.3690693333333333
$ws.Range("H12").Value = 1.107208
$ws.Range("I12").Value = 0.005124804994293743
$ws.Range("J12").Value = 0.005296051365434103
$ws.Range("M12").Value = 8.850733666666665
$ws.Range("N12").Value = 26.552201
$ws.Range("O12").Value = 0.03425243510433672
$ws.Range("P12").Value = 0.03471479367970558
$ws.Range("Q12").Value = 3.266534373867555
$ws.Range("R12").Value = 29.398809364808
$ws.Range("S12").Value = 0.0001755370504894271
$ws.Range("T12").Value = 0.0001838513304681679
$ws.Range("G13").Value = 0.3690693333333333
$ws.Range("H13").Value = 1.107208
$ws.Range("I13").Value = 0.005124804994293743
$ws.Range("J13").Value = 0.005296051365434103
$ws.Range("O13").Value = 0.06418224226698629
$ws.Range("P13").Value = 0.06504861016194469
$ws.Range("Q13").Value = 6.120834910521777
$ws.Range("R13").Value = 55.087514194696
$ws.Range("S13").Value = 0.0003289214757148223
$ws.Range("T13").Value = 0.0003445007806677578
$ws.Range("G14").Value = 0.3690693333333333
$ws.Range("H14").Value = 1.107208
$ws.Range("I14").Value = 0.005124804994293743
$ws.Range("J14").Value = 0.005296051365434103
$ws.Range("M14").Value = 94.82067633333332
$ws.Range("N14").Value = 284.462029
$ws.Range("O14").Value = 0.3669570438989389
$ws.Range("P14").Value = 0.3719104358409092
$ws.Range("Q14").Value = 34.9954038005591
$ws.Range("R14").Value = 314.958634205032
$ws.Range("S14").Value = 0.00188058329126455
$ws.Range("T14").Value = 0.001969656771554439
$ws.Range("G15").Value = 0.3690693333333333
$ws.Range("H15").Value = 1.107208
$ws.Range("I15").Value = 0.005124804994293743
$ws.Range("J15").Value = 0.005296051365434103
$ws.Range("M15").Value = 10.324604
$ws.Range("N15").Value = 20.649208
$ws.Range("O15").Value = 0.03995632924983982
$ws.Range("P15").Value = 0.02699712145781534
$ws.Range("Q15").Value = 3.810494715210667
$ws.Range("R15").Value = 22.862968291264
$ws.Range("S15").Value = 0.0002047683956932243
$ws.Range("T15").Value = 0.0001429781419594532
$ws.Range("G16").Value = 0.3690693333333333
$ws.Range("H16").Value = 1.107208
$ws.Range("I16").Value = 0.005124804994293743
$ws.Range("J16").Value = 0.005296051365434103
$ws.Range("M16").Value = 127.8166836666667
$ws.Range("N16").Value = 383.450051
$ws.Range("O16").Value = 0.4946519494798983
$ws.Range("P16").Value = 0.5013290388596253
$ws.Range("Q16").Value = 47.17321822973422
$ws.Range("R16").Value = 424.558964067608
$ws.Range("S16").Value = 0.002534994781131719
$ws.Range("T16").Value = 0.002655064340784285
$ws.Range("G17").Value = 6.985879
$ws.Range("H17").Value = 13.971758
$ws.Range("I17").Value = 0.0970041787687547
$ws.Range("J17").Value = 0.06683039504177611
$ws.Range("M17").Value = 8.850733666666665
$ws.Range("N17").Value = 26.552201
$ws.Range("O17").Value = 0.03425243510433672
$ws.Range("P17").Value = 0.03471479367970558
$ws.Range("Q17").Value = 61.83015445655965
$ws.Range("R17").Value = 370.9809267393579
$ws.Range("S17").Value = 0.003322629338126248
$ws.Range("T17").Value = 0.002320003375408476
$ws.Range("G18").Value = 6.985879
$ws.Range("H18").Value = 13.971758
$ws.Range("I18").Value = 0.0970041787687547
$ws.Range("J18").Value = 0.06683039504177611
$ws.Range("O18").Value = 0.06418224226698629
$ws.Range("P18").Value = 0.06504861016194469
$ws.Range("Q18").Value = 115.8573964346743
$ws.Range("R18").Value = 695.1443786080459
$ws.Range("S18").Value = 0.006225945702646262
$ws.Range("T18").Value = 0.004347224314041255
$ws.Range("G19").Value = 6.985879
$ws.Range("H19").Value = 13.971758
$ws.Range("I19").Value = 0.0970041787687547
$ws.Range("J19").Value = 0.06683039504177611
$ws.Range("M19").Value = 94.82067633333332
$ws.Range("N19").Value = 284.462029
$ws.Range("O19").Value = 0.3669570438989389
$ws.Range("P19").Value = 0.3719104358409092
$ws.Range("Q19").Value = 662.4057715628303
$ws.Range("R19").Value = 3974.434629376981
$ws.Range("S19").Value = 0.03559636668682643
$ws.Range("T19").Value = 0.02485492134740709
$ws.Range("G20").Value = 6.985879
$ws.Range("H20").Value = 13.971758
$ws.Range("I20").Value = 0.0970041787687547
$ws.Range("J20").Value = 0.06683039504177611
$ws.Range("M20").Value = 10.324604
$ws.Range("N20").Value = 20.649208
$ws.Range("O20").Value = 0.03995632924983982
$ws.Range("P20").Value = 0.02699712145781534
$ws.Range("Q20").Value = 72.12643426691601
$ws.Range("R20").Value = 288.505737067664
$ws.Range("S20").Value = 0.003875930905494684
$ws.Range("T20").Value = 0.00180422829201661
$ws.Range("G21").Value = 6.985879
$ws.Range("H21").Value = 13.971758
$ws.Range("I21").Value = 0.0970041787687547
$ws.Range("J21").Value = 0.06683039504177611
$ws.Range("M21").Value = 127.8166836666667
$ws.Range("N21").Value = 383.450051
$ws.Range("O21").Value = 0.4946519494798983
$ws.Range("P21").Value = 0.5013290388596253
$ws.Range("Q21").Value = 892.9118862766096
$ws.Range("R21").Value = 5357.471317659658
$ws.Range("S21").Value = 0.04798330613566108
$ws.Range("T21").Value = 0.03350401771290269
$ws.Range("G22").Value = 0.3102503333333334
$ws.Range("H22").Value = 0.930751
$ws.Range("I22").Value = 0.004308058985523854
$ws.Range("J22").Value = 0.004452013627456771
$ws.Range("M22").Value = 8.850733666666665
$ws.Range("N22").Value = 26.552201
$ws.Range("O22").Value = 0.03425243510433672
$ws.Range("P22").Value = 0.03471479367970558
$ws.Range("Q22").Value = 2.745943070327888
$ws.Range("R22").Value = 24.713487632951
$ws.Range("S22").Value = 0.0001475615108273105
$ws.Range("T22").Value = 0.0001545507345363994
$ws.Range("G23").Value = 0.3102503333333334
$ws.Range("H23").Value = 0.930751
$ws.Range("I23").Value = 0.004308058985523854
$ws.Range("J23").Value = 0.004452013627456771
$ws.Range("O23").Value = 0.06418224226698629
$ws.Range("P23").Value = 0.06504861016194469
$ws.Range("Q23").Value = 5.145350479587444
$ws.Range("R23").Value = 46.308154316287
$ws.Range("S23").Value = 0.0002765008855093593
$ws.Range("T23").Value = 0.0002895972988881007
$ws.Range("G24").Value = 0.3102503333333334
$ws.Range("H24").Value = 0.930751
$ws.Range("I24").Value = 0.004308058985523854
$ws.Range("J24").Value = 0.004452013627456771
$ws.Range("M24").Value = 94.82067633333332
$ws.Range("N24").Value = 284.462029
$ws.Range("O24").Value = 0.3669570438989389
$ws.Range("P24").Value = 0.3719104358409092
$ws.Range("Q24").Value = 29.41814643930878
$ws.Range("R24").Value = 264.763317953779
$ws.Range("S24").Value = 0.001580872590270095
$ws.Range("T24").Value = 0.001655750328557115
$ws.Range("G25").Value = 0.3102503333333334
$ws.Range("H25").Value = 0.930751
$ws.Range("I25").Value = 0.004308058985523854
$ws.Range("J25").Value = 0.004452013627456771
$ws.Range("M25").Value = 10.324604
$ws.Range("N25").Value = 20.649208
$ws.Range("O25").Value = 0.03995632924983982
$ws.Range("P25").Value = 0.02699712145781534
$ws.Range("Q25").Value = 3.203211832534667
$ws.Range("R25").Value = 19.219270995208
$ws.Range("S25").Value = 0.0001721342232533221
$ws.Range("T25").Value = 0.0001201915526322995
$ws.Range("G26").Value = 0.3102503333333334
$ws.Range("H26").Value = 0.930751
$ws.Range("I26").Value = 0.004308058985523854
$ws.Range("J26").Value = 0.004452013627456771
$ws.Range("M26").Value = 127.8166836666667
$ws.Range("N26").Value = 383.450051
$ws.Range("O26").Value = 0.4946519494798983
$ws.Range("P26").Value = 0.5013290388596253
$ws.Range("Q26").Value = 39.65516871314455
$ws.Range("R26").Value = 356.896518418301
$ws.Range("S26").Value = 0.002130989775663768
$ws.Range("T26").Value = 0.002231923712842857
